# Error Calculations and Plots
# Two rows ("RM 232" and "SC 92") were dropped from the missing-data
# table and the remaining rows shifted up; a few previously-blank cells
# were (re)populated / cleared as part of the re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 is "RM 232" - delete it entirely, shifting rows 27:35 up to 26:34.
$ws.Range("A26:F26").EntireRow.Delete()

# "SC 92" was originally row 28, now sits at row 27 after the shift above.
# Delete it too, shifting the remaining rows up to 26:33.
$ws.Range("A27:F27").EntireRow.Delete()

# Fix up the handful of cells whose values differ from a straight shift.
$ws.Range("C26").Value = 10.8      # SC 5 / column B value imputed
$ws.Range("C27").Value = ""        # SC 101 / column B value now missing
$ws.Range("E33").Value = -10.7     # SC 232 / column D value imputed
